$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge "Step 1. " + "Click this to fork" into a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Step 1. Click this to fork", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Step 1. Click this to fork", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: merge "Step 2. " + "Your account will appear in the popup.
# Click it" into a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Step 2. Your account will appear in the popup. Click it", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Step 2. Your account will appear in the popup. Click it", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3: the old "_GoBack" bookmark (at the end of the document, right
# after "Give your pull request a title ...") is removed -- it gets
# relocated to mark the "git remote add upstream" edit below.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# Change 4: fix "git add remote upstream" -> "git remote add upstream"
# (the word "add" was misplaced before "remote"). This is modeled as the
# interactive edit that produced it: cut the word "add " out from its
# original spot and paste it back in right after "remote ", which is
# what naturally splits the line into three runs with a collapsed
# "_GoBack" bookmark sitting at the point where the paste happened.
# ---------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("git add remote upstream")

# Normalize/re-seat the run first (strips any stale rsid bookkeeping from
# the original run so the pieces we are about to carve out of it come out
# clean).
$clean = $d.Range($target.Start, $target.End)
$clean.Find.Execute("git add remote upstream", $true, $false, $false, $false, $false,
                     $true, 1, $false, "git add remote upstream", 2) | Out-Null
$clean = $d.Range($target.Start, $target.End)

# Cut "add " out of "git add remote upstream".
$addWord = $d.Range($clean.Start, $clean.End)
$addWord.Find.Execute("add ") | Out-Null
$addWord.Cut()

# Paste "add " back in immediately after "remote ".
$afterRemote = $d.Range($clean.Start, $clean.End - 4)
$afterRemote.Find.Execute("remote ") | Out-Null
$afterRemote.Collapse(0)
$afterRemote.Paste()

# Drop a fresh collapsed "_GoBack" bookmark right where the paste landed,
# i.e. immediately before "upstream".
$beforeUpstream = $d.Range($clean.Start, $clean.End)
$beforeUpstream.Find.Execute("upstream") | Out-Null
$beforeUpstream.Collapse(1)
$d.Bookmarks.Add("_GoBack", $beforeUpstream) | Out-Null
